$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.604.57"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "'2.471.04"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'318.94"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").Value = "'92.30"
$ws.Range("E6").Value = "  +1.52%  "

$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").Value = "'0.0865"
$ws.Range("E10").Value = "  +9.23%  "

$ws.Range("D11").Value = "'33.01"
$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "'2.850.01"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").Value = "'2.474.91"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  +2.89%  "

$ws.Range("D18").Value = "'41.545.41"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "'6.45"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "'0.0₃0943"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").Value = "'70.79"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'11.28"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'239.96"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").Value = "'1.97"
$ws.Range("E25").Value = "  +3.65%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").Value = "'9.69"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").Value = "'36.49"
$ws.Range("E30").Value = "  +4.30%  "

$ws.Range("D31").Value = "'157.29"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").Value = "'0.0762"
$ws.Range("E34").Value = "  +0.97%  "

$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("D36").Value = "'17.27"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  +5.19%  "

$ws.Range("D38").Value = "'2.90"
$ws.Range("E38").Value = "  +1.69%  "

$ws.Range("E39").Value = "  +1.67%  "

$ws.Range("E40").Value = "  +1.98%  "

$ws.Range("E41").Value = "  +5.22%  "

$ws.Range("E42").Value = "  +0.58%  "

$ws.Range("D43").Value = "'1.988.17"
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D44").Value = "'0.0283"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("D45").Value = "'18.86"
$ws.Range("E45").Value = "  +1.23%  "

$ws.Range("D46").Value = "'2.96"
$ws.Range("E46").Value = "  +2.33%  "

$ws.Range("D47").Value = "'9.42"
$ws.Range("E47").Value = "  +5.08%  "

$ws.Range("D48").Value = "'2.703.94"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").Value = "'97.64"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").Value = "'75.78"
$ws.Range("E50").Value = "  +6.24%  "

$ws.Range("D51").Value = "'66.87"
$ws.Range("E51").Value = "  +0.27%  "
